$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 545, pushing the existing rows 545-557
# down to 547-559 (Excel copies formatting/row-style down automatically,
# matching the D-column date style already used throughout the table).
$ws.Rows("545:546").Insert()

# --- Row 545: new weekly price record ---
$ws.Range("A545").Value = 3
$ws.Range("B545").Value = "Femacal de La Calera"
$ws.Range("C545").Value = "Coquimbo"
$ws.Range("D545").Value = 44448
$ws.Range("E545").Value = 5
$ws.Range("F545").Value = "Fruta"
$ws.Range("G545").Value = 100102
$ws.Range("H545").Value = "Cítricos"
$ws.Range("I545").Value = 100102003
$ws.Range("J545").Value = "Limón"
$ws.Range("K545").Value = "Sin especificar"
$ws.Range("L545").Value = "1a amarillo"
$ws.Range("M545").Value = 365
$ws.Range("N545").Value = 3000
$ws.Range("O545").Value = 3500
$ws.Range("P545").Value = 3227
$ws.Range("Q545").Value = "$/malla 16 kilos"
$ws.Range("R545").Value = "Provincia de Quillota"
$ws.Range("S545").Value = 202
$ws.Range("T545").Value = 16

# --- Row 546: new weekly price record ---
$ws.Range("A546").Value = 3
$ws.Range("B546").Value = "Femacal de La Calera"
$ws.Range("C546").Value = "Coquimbo"
$ws.Range("D546").Value = 44448
$ws.Range("E546").Value = 5
$ws.Range("F546").Value = "Fruta"
$ws.Range("G546").Value = 100102
$ws.Range("H546").Value = "Cítricos"
$ws.Range("I546").Value = 100102003
$ws.Range("J546").Value = "Limón"
$ws.Range("K546").Value = "Sin especificar"
$ws.Range("L546").Value = "2a amarillo"
$ws.Range("M546").Value = 220
$ws.Range("N546").Value = 2400
$ws.Range("O546").Value = 2500
$ws.Range("P546").Value = 2445
$ws.Range("Q546").Value = "$/malla 16 kilos"
$ws.Range("R546").Value = "Provincia de Quillota"
$ws.Range("S546").Value = 153
$ws.Range("T546").Value = 16
